$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("04 Sep")

# Update data values in column E (the "Best" column) per the new results
$ws.Range("E12").Value = 336
$ws.Range("E16").Value = 65401
$ws.Range("E17").Value = 129269
$ws.Range("E21").Value = 140
$ws.Range("E22").Value = 278

# Update the active selection to reflect where the cursor was left
$ws.Range("E17").Select()
